$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings must be introduced in this exact first-use order so
# the generated sharedStrings.xml indices line up with the target workbook:
#   12 PR-004  (D5)
#   13 Summary (H1)
#   14 MCU hardware pin is not able to drive the MODE LED (H2)
#   15 MCU hardware pin is not able to drive the CRESET pin (H3)
#   16 Software defect results in MCU firmware lockup (H4)
#   17 FPGA is not able to drive MIDI THRU signals (H5)
#   18 MCU is not able to drive SCL/SDA signals. See PR-004 (H6)
#   19 PR-005 (D6)
#   20 Title (G1)
#   21 PCB - MODE LED (G2)
#   22 PCB - CRESET (G3)
#   23 SW - USB CDC Lockup (G4)
#   24 FW - MIDI THRU (G5)
#   25 PCB - I2C (G6)

# Row 5 (new problem report PR-004)
$ws.Range("D5").Value = "PR-004"

# Column H ("Summary") header + cells, top to bottom
$ws.Range("H1").Value = "Summary"
$ws.Range("H2").Value = "MCU hardware pin is not able to drive the MODE LED"
$ws.Range("H3").Value = "MCU hardware pin is not able to drive the CRESET pin"
$ws.Range("H4").Value = "Software defect results in MCU firmware lockup"
$ws.Range("H5").Value = "FPGA is not able to drive MIDI THRU signals"
$ws.Range("H6").Value = "MCU is not able to drive SCL/SDA signals. See PR-004"

# Row 6 (new problem report PR-005)
$ws.Range("D6").Value = "PR-005"

# Column G ("Title") header + cells, top to bottom
$ws.Range("G1").Value = "Title"
$ws.Range("G2").Value = "PCB - MODE LED"
$ws.Range("G3").Value = "PCB - CRESET"
$ws.Range("G4").Value = "SW - USB CDC Lockup"
$ws.Range("G5").Value = "FW - MIDI THRU"
$ws.Range("G6").Value = "PCB - I2C"

# --- Remaining cells that reuse already-existing shared strings ---
$ws.Range("E5").Value = "Open"
$ws.Range("F5").Value = "PCB"
$ws.Range("E6").Value = "Open"
$ws.Range("F6").Value = "PCB"

# --- Date cells: set the number format *before* assigning the value so the
# engine reuses the workbook's existing date style (s="1") instead of minting
# a brand-new cellXf/numFmt. ---
$dateFmt = $ws.Range("A2").NumberFormat

$ws.Range("B4").NumberFormat = $dateFmt
$ws.Range("B4").Value = (Get-Date -Year 2019 -Month 4 -Day 3 -Hour 0 -Minute 0 -Second 0)

$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("A5").Value = (Get-Date -Year 2019 -Month 4 -Day 3 -Hour 0 -Minute 0 -Second 0)

$ws.Range("A6").NumberFormat = $dateFmt
$ws.Range("A6").Value = (Get-Date -Year 2019 -Month 4 -Day 3 -Hour 0 -Minute 0 -Second 0)

# --- Column widths: new column G matches column F's (bestFit) width ---
$ws.Columns.Item(7).ColumnWidth = 11.5

# --- Selection, matching the saved workbook's cursor position ---
[void]$ws.Range("G7").Select()
